# Update Il21-Il21r.xlsx with new TPM-derived NATMI values.
# Target/sending cluster string table gains "Resolving-Mac" and the
# ligand-receptor table is recomputed + extended with the new
# "Resolving-Mac" sending-cluster rows (rows 6-9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) is unchanged. Columns A:T for data rows 2.. are:
#  A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
#  E..T = the NATMI-derived numeric metrics.
$rows = @(
  @("FAPs", "Il21", "Il21r", "ECs", 1, 0.3333333333333333, 0.068825, 0.206475, 0.6049060451053807, 0.6049060451053806, 3, 1, 0.9270496666666667, 2.781149, 0.05906853381764246, 0.05906853381764246, 0.06380419330833333, 0.574237739775, 0.03573091318180353, 0.03573091318180353),
  @("FAPs", "Il21", "Il21r", "FAPs", 1, 0.3333333333333333, 0.068825, 0.206475, 0.6049060451053807, 0.6049060451053806, 3, 1, 2.338564, 7.015692, 0.1490055513588678, 0.1490055513588678, 0.1609516673, 1.4485650057, 0.0901343587712394, 0.09013435877123936),
  @("FAPs", "Il21", "Il21r", "MuSCs", 1, 0.3333333333333333, 0.068825, 0.206475, 0.6049060451053807, 0.6049060451053806, 3, 1, 0.422343, 1.267029, 0.02691029690765713, 0.02691029690765712, 0.029067756975, 0.261609812775, 0.01627820127502243, 0.01627820127502242),
  @("FAPs", "Il21", "Il21r", "Resolving-Mac", 1, 0.3333333333333333, 0.068825, 0.206475, 0.6049060451053807, 0.6049060451053806, 3, 1, 12.006519, 36.019557, 0.7650156179158327, 0.7650156179158326, 0.8263486701749999, 7.437138031575, 0.4627625718773154, 0.4627625718773152),
  @("Resolving-Mac", "Il21", "Il21r", "ECs", 1, 0.3333333333333333, 0.044953, 0.134859, 0.3950939548946194, 0.3950939548946193, 3, 1, 0.9270496666666667, 2.781149, 0.05906853381764246, 0.05906853381764246, 0.04167366366566667, 0.3750629729910001, 0.02333762063583893, 0.02333762063583893),
  @("Resolving-Mac", "Il21", "Il21r", "FAPs", 1, 0.3333333333333333, 0.044953, 0.134859, 0.3950939548946194, 0.3950939548946193, 3, 1, 2.338564, 7.015692, 0.1490055513588678, 0.1490055513588678, 0.105125467492, 0.9461292074280001, 0.0588711925876284, 0.05887119258762839),
  @("Resolving-Mac", "Il21", "Il21r", "MuSCs", 1, 0.3333333333333333, 0.044953, 0.134859, 0.3950939548946194, 0.3950939548946193, 3, 1, 0.422343, 1.267029, 0.02691029690765713, 0.02691029690765712, 0.018985584879, 0.170870263911, 0.0106320956326347, 0.0106320956326347),
  @("Resolving-Mac", "Il21", "Il21r", "Resolving-Mac", 1, 0.3333333333333333, 0.044953, 0.134859, 0.3950939548946194, 0.3950939548946193, 3, 1, 12.006519, 36.019557, 0.7650156179158327, 0.7650156179158326, 0.539729048607, 4.857561437463001, 0.3022530460385174, 0.3022530460385173)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $rowNum = $i + 2
  $data = $rows[$i]
  $arr = New-Object "object[,]" 1,20
  for ($j = 0; $j -lt 20; $j++) { $arr[0,$j] = $data[$j] }
  $ws.Range("A$rowNum`:T$rowNum").Value = $arr
}
